# Cosenza.xlsx update: append new daily-data row (r=543, 2021-09-02) with
# its 7-day rolling average, plus placeholder date-only rows through
# r=571 (2021-09-30) on every data sheet, then restore each sheet's
# saved selection to the new last row (A543:D543) while leaving the
# workbook's active tab on sheet 3 ("Dimessi   Guariti"), matching the
# original file's activeTab.

$wb = $excel.ActiveWorkbook

# New case counts for 2021-09-02 (row 543), per sheet, in tab order:
# 1 Nuovi casi, 2 Deceduti, 3 Dimessi   Guariti, 4 Ricoveri, 5 Terapia
$newRowSerial = 44441
$newCValues = @{ 1 = 136; 2 = 2; 3 = 41; 4 = 42; 5 = 3 }
$sheetsWithStyledC = @(3, 5)

for ($s = 1; $s -le 5; $s++) {
    $ws = $wb.Worksheets.Item($s)

    # Row 543: date, new count, 7-day trailing average formula.
    $ws.Cells.Item(543, 1).Value = $newRowSerial
    $ws.Cells.Item(543, 3).Value = $newCValues[$s]

    if ($sheetsWithStyledC -contains $s) {
        $ws.Cells.Item(542, 3).Copy()
        $ws.Cells.Item(543, 3).PasteSpecial(-4122)
        $ws.Cells.Item(543, 3).Value = $newCValues[$s]
    }

    $ws.Cells.Item(542, 4).Copy()
    $ws.Cells.Item(543, 4).PasteSpecial(-4122)
    $ws.Cells.Item(543, 4).Formula = "=AVERAGE(C537:C543)"

    # Rows 544-571: forward-filled placeholder dates only (2021-09-03
    # through 2021-09-30), no case data yet.
    for ($i = 0; $i -lt 28; $i++) {
        $row = 544 + $i
        $serial = 44442 + $i
        $ws.Cells.Item($row, 1).Value = $serial
    }
}

# Move each sheet's saved selection to the new last row, cycling
# through every sheet and finishing on sheet 3 so that sheet stays the
# workbook's active tab (as in the original file).
foreach ($s in 1, 2, 4, 5, 3) {
    $ws = $wb.Worksheets.Item($s)
    $ws.Range("A543:D543").Select()
}
